$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C2").Value = 289
$ws.Range("C3").Value = 178839
$ws.Range("C4").Value = 168796
$ws.Range("C8").Value = 64.90000000000001
